$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("all")

# The "all" sheet had the 41-deg/48hr temp group (40.98083) duplicated
# right before the final 42.05973 group (rows 82-89 and 90-97 respectively).
# Remove the duplicate block (old rows 82-89); the 42.05973 rows shift up
# to become the new rows 82-89.
$ws.Rows("82:89").Delete()

# Reflect the new selection/scroll position shown in the saved workbook.
[void]$ws.Activate()
[void]$ws.Range("D83").Select()
$excel.ActiveWindow.ScrollRow = 79
$excel.ActiveWindow.ScrollColumn = 1
